{"js": "// Remove the \"Affiliate disclosure\" bullet paragraph from the body\n// (now communicated via the website sidebar instead of inline per-brief).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"Affiliate disclosure\") {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Affiliate disclosure\" bullet paragraph from the body\n# (now communicated via the website sidebar instead of inline per-brief).\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Affiliate disclosure\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $target.Range.Delete()\n}\n"}
